{"js": "// Add a \"Requisitos\" Heading2 section followed by a ListBullet paragraph\n// listing the weak prerequisite, appended at the end of the document body\n// (right after the Bibliografia text, before the section break).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// New Heading2 paragraph: \"Requisitos\"\nconst headingParagraph = lastParagraph.insertParagraph(\"Requisitos\", \"After\");\nheadingParagraph.styleBuiltIn = Word.BuiltInStyleName.heading2;\nawait context.sync();\n\n// New ListBullet paragraph with the prerequisite text, followed by a line\n// break (the trailing \"\\v\" becomes a <w:br/> inside the same run).\nconst bulletParagraph = headingParagraph.insertParagraph(\n  \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\\v\",\n  \"After\"\n);\nbulletParagraph.style = \"List Bullet\";\nawait context.sync();\n", "ps1": "# Add a \"Requisitos\" Heading 2 section followed by a List Bullet paragraph\n# listing the weak prerequisite, appended at the end of the document body\n# (right after the Bibliografia text, before the section break).\n\n$d = $word.ActiveDocument\n\n# --- New Heading 2 paragraph: \"Requisitos\" ---------------------------------\n$lastParagraph = $d.Paragraphs.Last\n$tailRange = $lastParagraph.Range\n$tailRange.Collapse(0)            # wdCollapseEnd\n$tailRange.InsertParagraphAfter()\n\n$headingParagraph = $d.Paragraphs.Last\n$headingParagraph.Range.Text = \"Requisitos\"\n$headingParagraph.Style = \"Heading 2\"\n\n# --- New List Bullet paragraph with the prerequisite + trailing line break -\n$headingTail = $headingParagraph.Range\n$headingTail.Collapse(0)          # wdCollapseEnd\n$headingTail.InsertParagraphAfter()\n\n$bulletParagraph = $d.Paragraphs.Last\n$bulletParagraph.Range.Text = \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\"\n$bulletParagraph.Style = \"List Bullet\"\n\n$bulletTail = $bulletParagraph.Range\n$bulletTail.Collapse(0)           # wdCollapseEnd\n$bulletTail.InsertBefore([char]11)  # manual line break (<w:br/>)\n"}
